$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.656.58'
$ws.Range('E2').Value = '  +3.13%  '
$ws.Range('D3').Value = '3.296.40'
$ws.Range('E3').Value = '  +0.50%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.ClearFormats()
$ws.Range('E4').Value = '  +0.03%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '580.78'
$c.ClearFormats()
$ws.Range('E5').Value = '  +0.85%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '179.27'
$c.ClearFormats()
$ws.Range('E6').Value = '  -1.27%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  +3.17%  '
$ws.Range('D9').Value = '3.288.52'
$ws.Range('E9').Value = '  +0.39%  '
$ws.Range('E10').Value = '  +1.26%  '
$ws.Range('E11').Value = '  +1.57%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '45.85'
$c.ClearFormats()
$ws.Range('E12').Value = '  -0.79%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.0000272'
$c.ClearFormats()
$ws.Range('E13').Value = '  +3.76%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '694.09'
$c.ClearFormats()
$ws.Range('E14').Value = '  +13.80%  '
$ws.Range('D15').Value = '3.823.24'
$ws.Range('E15').Value = '  +0.51%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '8.41'
$c.ClearFormats()
$ws.Range('E16').Value = '  +0.56%  '
$ws.Range('D17').Value = '67.775.08'
$ws.Range('E17').Value = '  +3.25%  '
$ws.Range('E18').Value = '  +1.52%  '
$ws.Range('D19').Value = '3.292.50'
$ws.Range('E19').Value = '  +0.32%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '17.50'
$c.ClearFormats()
$ws.Range('E20').Value = '  -0.66%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '10.84'
$c.ClearFormats()
$ws.Range('E21').Value = '  -0.08%  '
$ws.Range('E22').Value = '  +1.54%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '17.01'
$c.ClearFormats()
$ws.Range('E23').Value = '  -6.79%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '5.22'
$c.ClearFormats()
$ws.Range('E24').Value = '  +6.29%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '98.32'
$c.ClearFormats()
$ws.Range('E25').Value = '  +0.42%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '3.96'
$c.ClearFormats()
$ws.Range('E26').Value = '  +0.54%  '
$ws.Range('E27').Value = '  +1.66%  '
$ws.Range('E28').Value = '  -0.01%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '33.18'
$c.ClearFormats()
$ws.Range('E29').Value = '  +8.34%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '8.50'
$c.ClearFormats()
$ws.Range('E30').Value = '  +1.87%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '6.79'
$c.ClearFormats()
$ws.Range('E31').Value = '  +5.49%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '586.36'
$c.ClearFormats()
$ws.Range('E32').Value = '  +7.54%  '
$ws.Range('E33').Value = '  +1.08%  '
$ws.Range('D34').Value = '3.879.56'
$ws.Range('E34').Value = '  +2.22%  '
$ws.Range('E35').Value = '  +1.82%  '
$ws.Range('E36').Value = '  +0.24%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '3.41'
$c.ClearFormats()
$ws.Range('E37').Value = '  -8.17%  '
$ws.Range('E38').Value = '  -1.20%  '
$ws.Range('E39').Value = '  +2.33%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '3.21'
$c.ClearFormats()
$ws.Range('E40').Value = '  +2.44%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '2.64'
$c.ClearFormats()
$ws.Range('E41').Value = '  +3.00%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '32.39'
$c.ClearFormats()
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('B43').Value = 'PEPE'
$ws.Range('C43').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D43').Value = '0.0₃0689'
$ws.Range('E43').Value = '  +1.71%  '
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '3.39'
$c.ClearFormats()
$ws.Range('E44').Value = '  -0.37%  '
$ws.Range('E45').Value = '  +1.01%  '
$ws.Range('E46').Value = '  +2.50%  '
$ws.Range('E47').Value = '  +2.42%  '
$ws.Range('E48').Value = '  +10.41%  '
$ws.Range('E49').Value = '  +0.46%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '2.54'
$c.ClearFormats()
$ws.Range('E50').Value = '  +1.92%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '128.80'
$c.ClearFormats()
$ws.Range('E51').Value = '  +0.81%  '
